$d = $word.ActiveDocument

# First paragraph: the hidden **ID__...__ID** marker paragraph.
$p1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right) with 5-twip spacing,
# and no visible line (matches the "last minute update" styling tweak).
$b = $p1.Borders
$b.DistanceFromTop = 5
$b.DistanceFromLeft = 5
$b.DistanceFromBottom = 5
$b.DistanceFromRight = 5

# Bump the left indent from 120 -> 225 twips (11.25 pt).
$p1.LeftIndent = 11.25

# Replace the marker text and drop the trailing lone-space run, collapsing
# the paragraph down to a single run.
$r = $p1.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "**ID__AFFARS_SUBPART_5306_2__ID**"
